$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Rows 25-30: set Assigned = Richard, Date Done = 5/30/2021
foreach ($r in 25..30) {
    $ws.Range("D$r").Value = "Richard"
    $ws.Range("E$r").Value = [DateTime]"2021-05-30"
}

# Row 31: Assigned = Richard, Date Done = 5/30/2021 (cells already exist, just empty)
$ws.Range("D31").Value = "Richard"
$ws.Range("E31").Value = [DateTime]"2021-05-30"

# Row 32: Assigned = Alex, Date Done = 5/30/2021
$ws.Range("D32").Value = "Alex"
$ws.Range("E32").Value = [DateTime]"2021-05-30"

# Update sheet view: scroll position and selection
$ws.Range("E32").Select()
$excel.ActiveWindow.ScrollRow = 22
